# Doing Updates for Financials
# Shifts a new reporting period into column D (pushing older periods right,
# dropping the oldest period from column J), and marks now-missing data
# points as "NA", for each data row of the UQM yearly financials sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("UQM")

# Each entry: row number, then values for columns D,E,F,G,H,I,J (in order).
# A value of "NA" represents the literal text "NA" used elsewhere in the sheet.
$rowUpdates = @(
    @{ Row = 8;   Values = @(14200, 7800, 4100, 5300, 4000, 7000, 7200) },
    @{ Row = 9;   Values = @(10800, 4700, 10100, 4000, 2900, 4400, 5000) },
    @{ Row = 10;  Values = @(3400, 3100, -5900, 1300, 1100, 2700, 2100) },
    @{ Row = 12;  Values = @(2500, 2000, 2400, 3500, 4500, 4900, 5000) },
    @{ Row = 14;  Values = @("NA", 0, 0, -600, 0, -900, 5000) },
    @{ Row = 17;  Values = @(20600, 13100, 17200, 12300, 10900, 9800, 17900) },
    @{ Row = 18;  Values = @(-6400, -5300, -13000, -7000, -6900, -2800, -10700) },
    @{ Row = 20;  Values = @("NA", 600, 0, 0, 900, 0, 0) },
    @{ Row = 21;  Values = @(400, -4200, -12100, -5900, -4900, -1500, "NA") },
    @{ Row = 22;  Values = @("NA", 100, "NA", "NA", "NA", "NA", "NA") },
    @{ Row = 23;  Values = @("NA", -4800, -13000, -6900, -6000, -2800, -10700) },
    @{ Row = 26;  Values = @("NA", -4800, -13000, -6900, -6000, -2800, -10700) },
    @{ Row = 27;  Values = @("NA", -4800, -13000, -6900, -6000, -2800, -10700) },
    @{ Row = 32;  Values = @("NA", -600, 0, 0, -900, 0, 0) },
    @{ Row = 33;  Values = @("NA", -4800, -13000, -6900, -6000, -2800, -10700) },
    @{ Row = 35;  Values = @("NA", -4800, -13000, -6900, -6000, -2800, -10700) },
    @{ Row = 81;  Values = @("NA", -4800, -13000, -6900, -6000, -2800, -10700) },
    @{ Row = 83;  Values = @(400, 500, 1000, 1100, 1100, 1300, "NA") },
    @{ Row = 91;  Values = @(-100, 0, -100, -600, -400, -600, -2100) },
    @{ Row = 94;  Values = @(1300, -100, -200, -200, 1300, 100, "NA") },
    @{ Row = 100; Values = @(8300, 0, 5700, -100, 5000, 0, "NA") },
    @{ Row = 101; Values = @(0, 0, 0, 0, 0, 0, "NA") }
)

$cols = @("D", "E", "F", "G", "H", "I", "J")

foreach ($update in $rowUpdates) {
    $r = $update.Row
    $vals = $update.Values
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range("$($cols[$i])$r").Value = $vals[$i]
    }
}
